# "Generate Report for Archive"
#
# The handback status moved on from "Ready for handoff" to "In Translation"
# for the e2e\19175915-1708-4e2e-af2b-7cfe3d5a7cba.md asset, so every place
# that showed the old status text needs to show the new one. The same
# shared string backed the "Overview" sheet's per-language status columns
# (zh-cn / de-de) as well as the Status column on each language detail
# sheet, so all of them are updated together.
#
# Because the new status text is shorter than the old one, Excel's column
# autosizing also shrank the affected "Status" columns when the sheet was
# regenerated for the archive, so we bring those columns back in to the
# same (narrower) width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: per-language status cells (columns E = zh-cn, F = de-de)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-language detail sheets: Status column (column C)
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# The Status columns re-autofit narrower now that "In Translation" is
# shorter than "Ready for handoff". (The engine's ColumnWidth setter
# quantizes to the nearest 1/6 of a character, so 12.5 is the input that
# lands closest to the archived 13.4101848602295 width.)
$newStatusColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth
$zhcn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
$dede.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
